$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 2 de Septiembre de 2020 a las 22:03"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 6285556
$ws.Range("C4").Value = 27985
$ws.Range("D4").Value = 3527920
$ws.Range("E4").Value = 2568033
$ws.Range("G4").Value = 703
$ws.Range("H4").Value = 189603

# Row 9 - Sudafrica
$ws.Range("B9").Value = 630595
$ws.Range("C9").Value = 2336
$ws.Range("D9").Value = 553456
$ws.Range("E9").Value = 62750
$ws.Range("G9").Value = 126
$ws.Range("H9").Value = 14389

# Row 23 - Alemania
$ws.Range("B23").Value = 247152
$ws.Range("C23").Value = 1151
$ws.Range("E23").Value = 15963
$ws.Range("G23").Value = 8
$ws.Range("H23").Value = 9389

# Row 29 - Israel
$ws.Range("B29").Value = 121464
$ws.Range("C29").Value = 2926
$ws.Range("D29").Value = 97234
$ws.Range("E29").Value = 23261

# Row 97 - Guayana Francesa
$ws.Range("B97").Value = 9209
$ws.Range("C97").Value = 55
$ws.Range("D97").Value = 8739
$ws.Range("E97").Value = 409
$ws.Range("G97").Value = 1
$ws.Range("H97").Value = 61
